$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.856.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.51"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5018"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2573"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06416"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.65"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07699"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.244"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.57"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.861.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7946"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.50"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.871.86"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.329"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.954"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.977"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.928"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +11.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.28"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.712"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05007"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.263"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.177"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.535"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.363"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.174.58"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8949"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.603"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5622"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8068"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.48"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.773.16"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4517"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.80"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05060"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.01%  "
